$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 ("stethoscope_5_realHeart_" / "...channel_6.wav") was missing its
# raw S1/S2 counts and the per-category tally columns that every other data
# row (e.g. row 25) already has. Fill them in the same way Excel would if a
# user typed the numbers in: the columns already carry style 1 via the
# sheet's <cols> defaults, so a plain value write picks up the right style.

$ws.Range("C26").Value2 = 113
$ws.Range("D26").Value2 = 112
$ws.Range("E26").Value2 = 0

$ws.Range("L26").Value2 = 0
$ws.Range("M26").Value2 = 0
$ws.Range("N26").Value2 = 0
$ws.Range("O26").Value2 = 0
$ws.Range("P26").Value2 = 0
$ws.Range("Q26").Value2 = 0
$ws.Range("R26").Value2 = 0
$ws.Range("S26").Value2 = 0
$ws.Range("T26").Value2 = 0
